$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9873417721518988
$ws.Range("C2").Value = 0.7722772277227723
$ws.Range("D2").Value = 0.8666666666666666

$ws.Range("B3").Value = 0.8130081300813008
$ws.Range("C3").Value = 0.9900990099009901
$ws.Range("D3").Value = 0.8928571428571429

$ws.Range("B4").Value = 0.8811881188118812
$ws.Range("C4").Value = 0.8811881188118812
$ws.Range("D4").Value = 0.8811881188118812
$ws.Range("E4").Value = 0.8811881188118812

$ws.Range("B5").Value = 0.9001749511165997
$ws.Range("C5").Value = 0.8811881188118812
$ws.Range("D5").Value = 0.8797619047619047

$ws.Range("B6").Value = 0.9001749511165998
$ws.Range("C6").Value = 0.8811881188118812
$ws.Range("D6").Value = 0.8797619047619047
